# Update "paises.xlsx" (Pais sheet) with the 15-Abril-2020 02:22 data refresh.
# This mirrors the upstream scraper commit: country case numbers are refreshed
# and a handful of countries swap rank (since the sheet is kept sorted by
# "Casos totales" descending), plus the "last updated" timestamp changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    if ($country -ne $null) { $ws.Cells.Item($row, 1).Value = $country }
    if ($b -ne $null) { $ws.Cells.Item($row, 2).Value = $b }
    if ($c -ne $null) { $ws.Cells.Item($row, 3).Value = $c }
    if ($d -ne $null) { $ws.Cells.Item($row, 4).Value = $d }
    if ($e -ne $null) { $ws.Cells.Item($row, 5).Value = $e }
    if ($f -ne $null) { $ws.Cells.Item($row, 6).Value = $f }
    if ($g -ne $null) { $ws.Cells.Item($row, 7).Value = $g }
    if ($h -ne $null) { $ws.Cells.Item($row, 8).Value = $h }
}

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 02:22"

# --- Straight numeric refreshes (no reordering) ----------------------------
# Estados Unidos (row 4)
Set-Row 4 $null 613624 26683 38721 548887 13473 2376 26016

# Canada (row 15)
Set-Row 15 $null 27063 1383 8235 17925 557 123 903

# Brasil (row 17) - only Casos activos / Recuperados move
Set-Row 17 $null $null $null 14026 9704 $null $null $null

# Panama (row 44)
Set-Row 44 $null 3574 102 72 3407 106 1 95

# San Marino (row 105) - only Casos totales / Nuevos casos / Recuperados move
Set-Row 105 $null 372 1 $null 283 $null $null $null

# --- Re-ranked block: Guinea-Bisau jumps ahead of Guinea Ecuatorial / Haiti /
#     Puerto Rico (rows 157-160) ------------------------------------------
Set-Row 157 "Guinea-Bisau" 43 5 0 43 0 0 0
Set-Row 158 "Guinea Ecuatorial" 41 20 4 37 0 0 0
Set-Row 159 "Haiti" 40 0 0 37 0 0 3
Set-Row 160 "Puerto Rico" 39 0 1 36 0 0 2

# --- Re-ranked pair: Granada overtakes San Cristobal y Nieves (rows 186-187)
Set-Row 186 "Granada" 14 0 0 14 2 0 0
Set-Row 187 "San Cristobal y Nieves" 14 2 0 14 0 0 0

# --- Re-ranked block: Sierra Leona overtakes Seychelles, and Islas Malvinas
#     overtakes Montserrat (rows 191-194) -----------------------------------
Set-Row 191 "Sierra Leona" 11 1 0 11 0 0 0
Set-Row 192 "Seychelles" 11 0 0 11 0 0 0
Set-Row 193 "Islas Malvinas" 11 6 1 10 0 0 0
Set-Row 194 "Montserrat" 11 0 1 10 1 0 0
